# Auto-generated edit script: updates market-price derived columns (H-N)
# on sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR per the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2516.25
$ws.Range("I19").Value = 1350
$ws.Range("J19").Value = 2749.5
$ws.Range("K19").Value = 1350
$ws.Range("L19").Value = 2749.5
$ws.Range("M19").Value = -1175
$ws.Range("N19").Value = -3099.5
$ws.Range("H32").Value = 6159.8
$ws.Range("I32").Value = 7666.3335
$ws.Range("J32").Value = 3900
$ws.Range("K32").Value = 7666.3335
$ws.Range("L32").Value = 3900
$ws.Range("M32").Value = -7340.3335
$ws.Range("N32").Value = -4552
$ws.Range("H70").Value = 7143.0713
$ws.Range("I70").Value = 7666.6665
$ws.Range("J70").Value = 4001.5
$ws.Range("K70").Value = 22999.9995
$ws.Range("L70").Value = 12004.5
$ws.Range("M70").Value = -22729.9995
$ws.Range("N70").Value = -12544.5
$ws.Range("H73").Value = 7143.0713
$ws.Range("I73").Value = 7666.6665
$ws.Range("J73").Value = 4001.5
$ws.Range("K73").Value = 22999.9995
$ws.Range("L73").Value = 12004.5
$ws.Range("M73").Value = -22063.9995
$ws.Range("N73").Value = -13876.5
$ws.Range("H94").Value = 3162.25
$ws.Range("I94").Value = 3399.7144
$ws.Range("K94").Value = 3399.7144
$ws.Range("M94").Value = -2948.7144
$ws.Range("H97").Value = 2948.4285
$ws.Range("J97").Value = 2036.5
$ws.Range("L97").Value = 6109.5
$ws.Range("N97").Value = -7101.5
$ws.Range("H101").Value = 343.42856
$ws.Range("I101").Value = 404
$ws.Range("K101").Value = 1212
$ws.Range("M101").Value = 410
$ws.Range("H107").Value = 2605.375
$ws.Range("I107").Value = 2727.5715
$ws.Range("K107").Value = 2727.5715
$ws.Range("M107").Value = -807.5715
$ws.Range("H127").Value = 2871.2144
$ws.Range("I127").Value = 2895.75
$ws.Range("J127").Value = 2724
$ws.Range("K127").Value = 8687.25
$ws.Range("L127").Value = 8172
$ws.Range("M127").Value = -3727.25
$ws.Range("N127").Value = -18092
$ws.Range("H135").Value = 1331.5555
$ws.Range("I135").Value = 1331.5555
$ws.Range("K135").Value = 11983.9995
$ws.Range("M135").Value = -9448.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 348338.53
$ws.Range("I32").Value = 1510.381
$ws.Range("K32").Value = 1510.381
$ws.Range("M32").Value = -1223.381
$ws.Range("H45").Value = 1395.9286
$ws.Range("I45").Value = 1407.9
$ws.Range("J45").Value = 1366
$ws.Range("K45").Value = 1407.9
$ws.Range("L45").Value = 1366
$ws.Range("M45").Value = -1030.9
$ws.Range("N45").Value = -2120
$ws.Range("H56").Value = 39245
$ws.Range("J56").Value = 39245
$ws.Range("L56").Value = 39245
$ws.Range("N56").Value = -40729
$ws.Range("H110").Value = 1514.7931
$ws.Range("I110").Value = 1129.2307
$ws.Range("K110").Value = 1129.2307
$ws.Range("M110").Value = 915.7692999999999
$ws.Range("H122").Value = 2724.8462
$ws.Range("J122").Value = 5000
$ws.Range("L122").Value = 15000
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 2809.2666
$ws.Range("I132").Value = 2761.75
$ws.Range("J132").Value = 2999.3333
$ws.Range("K132").Value = 8285.25
$ws.Range("L132").Value = 8997.999899999999
$ws.Range("M132").Value = -5755.25
$ws.Range("N132").Value = -14057.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2666.6667
$ws.Range("J94").Value = 1500
$ws.Range("L94").Value = 1500
$ws.Range("N94").Value = -2402
$ws.Range("H99").Value = 1157.0869
$ws.Range("I99").Value = 1132.409
$ws.Range("K99").Value = 1132.409
$ws.Range("M99").Value = 365.5909999999999
$ws.Range("H107").Value = 2092.2666
$ws.Range("I107").Value = 1865.3334
$ws.Range("K107").Value = 1865.3334
$ws.Range("M107").Value = 54.66660000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 54966.26
$ws.Range("I22").Value = 85524.664
$ws.Range("K22").Value = 85524.664
$ws.Range("M22").Value = -85174.664
$ws.Range("H58").Value = 2409.476
$ws.Range("I58").Value = 2090.2
$ws.Range("J58").Value = 3207.6667
$ws.Range("K58").Value = 2090.2
$ws.Range("L58").Value = 3207.6667
$ws.Range("M58").Value = -1887.2
$ws.Range("N58").Value = -3613.6667
$ws.Range("H62").Value = 201499.5
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376
$ws.Range("H65").Value = 201499.5
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880
$ws.Range("H86").Value = 8263.5
$ws.Range("I86").Value = 8234.333000000001
$ws.Range("K86").Value = 8234.333000000001
$ws.Range("M86").Value = -7111.333000000001
$ws.Range("H89").Value = 8263.5
$ws.Range("I89").Value = 8234.333000000001
$ws.Range("K89").Value = 41171.665
$ws.Range("M89").Value = -35555.665
$ws.Range("H94").Value = 2215.6667
$ws.Range("J94").Value = 2242.75
$ws.Range("L94").Value = 2242.75
$ws.Range("N94").Value = -3144.75
$ws.Range("H99").Value = 2877.5293
$ws.Range("I99").Value = 1925.5454
$ws.Range("K99").Value = 1925.5454
$ws.Range("M99").Value = -427.5454
$ws.Range("H126").Value = 2877.5293
$ws.Range("I126").Value = 1925.5454
$ws.Range("K126").Value = 5776.6362
$ws.Range("M126").Value = -3306.6362
$ws.Range("H132").Value = 6129
$ws.Range("I132").Value = 5491.1
$ws.Range("K132").Value = 16473.3
$ws.Range("M132").Value = -13943.3
$ws.Range("H134").Value = 1828.826
$ws.Range("I134").Value = 1506.1765
$ws.Range("K134").Value = 4518.529500000001
$ws.Range("M134").Value = -1983.529500000001
$ws.Range("H136").Value = 2409.476
$ws.Range("I136").Value = 2090.2
$ws.Range("J136").Value = 3207.6667
$ws.Range("K136").Value = 6270.599999999999
$ws.Range("L136").Value = 9623.000100000001
$ws.Range("M136").Value = -3720.599999999999
$ws.Range("N136").Value = -14723.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 683506.6
$ws.Range("I4").Value = 784637.9399999999
$ws.Range("K4").Value = 2353913.82
$ws.Range("M4").Value = -2353801.82
$ws.Range("H57").Value = 8873.75
$ws.Range("I57").Value = 8665.333000000001
$ws.Range("K57").Value = 25995.999
$ws.Range("M57").Value = -25436.999
$ws.Range("H93").Value = 1624.75
$ws.Range("I93").Value = 1624.75
$ws.Range("K93").Value = 4874.25
$ws.Range("M93").Value = -3002.25
$ws.Range("H94").Value = 3160
$ws.Range("J94").Value = 5000
$ws.Range("L94").Value = 15000
$ws.Range("N94").Value = -16352
$ws.Range("H96").Value = 10000
$ws.Range("J96").Value = 10000
$ws.Range("L96").Value = 30000
$ws.Range("N96").Value = -34118
$ws.Range("H107").Value = 91221.82000000001
$ws.Range("I107").Value = 199
$ws.Range("J107").Value = 111449.11
$ws.Range("K107").Value = 597
$ws.Range("L107").Value = 334347.33
$ws.Range("M107").Value = 1323
$ws.Range("N107").Value = -338187.33
$ws.Range("H108").Value = 8356.861999999999
$ws.Range("I108").Value = 2058.1667
$ws.Range("K108").Value = 6174.500100000001
$ws.Range("M108").Value = -3294.500100000001
$ws.Range("H128").Value = 953289.8
$ws.Range("I128").Value = 953289.8
$ws.Range("K128").Value = 2859869.4
$ws.Range("M128").Value = -2854889.4
$ws.Range("H139").Value = 2904.1667
$ws.Range("I139").Value = 1048
$ws.Range("K139").Value = 3144
$ws.Range("M139").Value = 1996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5210.222
$ws.Range("I132").Value = 4146.3335
$ws.Range("K132").Value = 12439.0005
$ws.Range("M132").Value = -9909.000499999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4263.1
$ws.Range("J7").Value = 3996.5
$ws.Range("L7").Value = 3996.5
$ws.Range("N7").Value = -4220.5
$ws.Range("H22").Value = 1191.4615
$ws.Range("I22").Value = 1155.5714
$ws.Range("J22").Value = 1233.3334
$ws.Range("K22").Value = 1155.5714
$ws.Range("L22").Value = 1233.3334
$ws.Range("M22").Value = -860.5714
$ws.Range("N22").Value = -1823.3334
$ws.Range("H27").Value = 1191.4615
$ws.Range("I27").Value = 1155.5714
$ws.Range("J27").Value = 1233.3334
$ws.Range("K27").Value = 1155.5714
$ws.Range("L27").Value = 1233.3334
$ws.Range("M27").Value = -1048.5714
$ws.Range("N27").Value = -1447.3334
$ws.Range("H40").Value = 2117.4666
$ws.Range("I40").Value = 2147.2856
$ws.Range("K40").Value = 2147.2856
$ws.Range("M40").Value = -2011.2856
$ws.Range("H126").Value = 4263.1
$ws.Range("J126").Value = 3996.5
$ws.Range("L126").Value = 11989.5
$ws.Range("N126").Value = -16929.5
$ws.Range("H132").Value = 133008.75
$ws.Range("I132").Value = 146295.72
$ws.Range("K132").Value = 438887.16
$ws.Range("M132").Value = -436357.16

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1532.381
$ws.Range("I122").Value = 1383.8235
$ws.Range("K122").Value = 4151.470499999999
$ws.Range("M122").Value = -1701.470499999999

Write-Host "Updated 228 cells across 8 sheets"
